# Generate Report for Handoff
#
# The handoff run id changed from 50499540-ca90-4829-844d-ebf5f83622b0 to
# 4c5f95a2-556f-46b2-84cd-e178c29d873d, and its content hash changed from
# 501959762e364b43bfb8ddbbf57690317b1083d8 to
# 2249aff29b95a65b98cd5c0145f982c65d6d25d3. Update the file-name cells and
# their timestamps on the Overview/zh-cn/de-de sheets, keeping each
# hyperlink's displayed text synchronized with its cell value.

$wb = $excel.ActiveWorkbook

$newId = "4c5f95a2-556f-46b2-84cd-e178c29d873d"
$newHash = "2249aff29b95a65b98cd5c0145f982c65d6d25d3"

function Set-CellAndHyperlink($ws, $addr, $newValue) {
    $ws.Range($addr).Value = $newValue
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $newValue
        }
    }
}

# --- Overview sheet: A2 = handoff file name, D2 = latest handoff date ---
$wsOverview = $wb.Worksheets.Item("Overview")
Set-CellAndHyperlink $wsOverview "`$A`$2" "$newId.md"
$wsOverview.Range("D2").Value = "2016-39-11 08:39:32"

# --- zh-cn sheet: A2 = source file name, D2 = handoff target, E2 = datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Set-CellAndHyperlink $wsZh "`$A`$2" "$newId.md"
Set-CellAndHyperlink $wsZh "`$D`$2" "$newId.$newHash.zh-cn.xlf"
$wsZh.Range("E2").Value = "2016-03-11 08:39:29"

# --- de-de sheet: A2 = source file name, D2 = handoff target, E2 = datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
Set-CellAndHyperlink $wsDe "`$A`$2" "$newId.md"
Set-CellAndHyperlink $wsDe "`$D`$2" "$newId.$newHash.de-de.xlf"
$wsDe.Range("E2").Value = "2016-03-11 08:39:32"
